# Apply the edits described by the diff:
#  1. Rename the worksheet from "GossF-HW40.xpc" to "GossF".
#  2. Append a new data row (row 16) re-using the existing "HexGrid-60degTilt5degRes"
#     label (same text as row 15) together with its computed averaged-intensity values.
#  3. Copy the formatting (bold/border/alignment style) from A15 into A16, matching
#     the style used by the rest of column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet.
$ws.Name = "GossF"

# 2. Populate the new row of data (row 16).
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "HexGrid-60degTilt5degRes"
$ws.Cells.Item(16, 3).Value = 1.001031114751944
$ws.Cells.Item(16, 4).Value = 0.9872471667522974
$ws.Cells.Item(16, 5).Value = 1.000508245343851
$ws.Cells.Item(16, 6).Value = 1.001031114751944
$ws.Cells.Item(16, 7).Value = 0.9914525581965761
$ws.Cells.Item(16, 8).Value = 1.002854316891267
$ws.Cells.Item(16, 9).Value = 0.9994117647058823
$ws.Cells.Item(16, 10).Value = 0.9872471667522974
$ws.Cells.Item(16, 11).Value = 0.9938777060480739
$ws.Cells.Item(16, 12).Value = 0.9974544104000088
$ws.Cells.Item(16, 13).Value = 0.9970841944403026

# 3. Match column-A styling used by the preceding rows (bold, bordered, centered/top).
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$excel.CutCopyMode = 0
